$wb = $excel.ActiveWorkbook

# --- Insert "detail_qa" sheet before "ticket_qa" ---
$wsTicketQa = $wb.Worksheets.Item("ticket_qa")
$detailQa = $wb.Worksheets.Add($wsTicketQa)
$detailQa.Name = "detail_qa"

# Header row
$detailQa.Range("B2").Value = "Name"
$detailQa.Range("C2").Value = "Data Type"
$detailQa.Range("D2").Value = "Primary Identifier"
$detailQa.Range("E2").Value = "Description"
$detailQa.Range("B2:D2").Font.Bold = $true
$detailQa.Range("B2:D2").Borders.LineStyle = 1
$detailQa.Range("B2:D2").Borders.Weight = 2
$detailQa.Range("E2").Font.Bold = $true
$detailQa.Range("E2").Borders.LineStyle = 1
$detailQa.Range("E2").Borders.Weight = 2
$detailQa.Range("E2").VerticalAlignment = -4108

# Row 3: dqa_id / Serial / TRUE / Llave primaria
$detailQa.Range("B3").Value = "dqa_id"
$detailQa.Range("C3").Value = "Serial"
$detailQa.Range("D3").Value = "TRUE"
$detailQa.Range("E3").Value = "Llave primaria"
$detailQa.Range("B3:D3").Borders.LineStyle = 1
$detailQa.Range("B3:D3").Borders.Weight = 2
$detailQa.Range("E3").Borders.LineStyle = 1
$detailQa.Range("E3").Borders.Weight = 2
$detailQa.Range("E3").VerticalAlignment = -4108

# Row 4: dtk_id / Integer / FALSE / id del detalle del ticket
$detailQa.Range("B4").Value = "dtk_id"
$detailQa.Range("C4").Value = "Integer"
$detailQa.Range("D4").Value = "FALSE"
$detailQa.Range("E4").Value = "id del detalle del ticket"
$detailQa.Range("B4:E4").Borders.LineStyle = 1
$detailQa.Range("B4:E4").Borders.Weight = 2
$detailQa.Range("B4:E4").VerticalAlignment = -4108

# Row 5: tqa_id / Integer / FALSE / id del ticket de QA
$detailQa.Range("B5").Value = "tqa_id"
$detailQa.Range("C5").Value = "Integer"
$detailQa.Range("D5").Value = "FALSE"
$detailQa.Range("E5").Value = "id del ticket de QA"
$detailQa.Range("B5:E5").Borders.LineStyle = 1
$detailQa.Range("B5:E5").Borders.Weight = 2

$detailQa.Columns.Item(4).ColumnWidth = 16.85546875
$detailQa.Columns.Item(5).ColumnWidth = 21.7109375

Write-Host "detail_qa created"

# --- Insert "detail_prod" sheet before "ticket_prod" ---
$wsTicketProd = $wb.Worksheets.Item("ticket_prod")
$detailProd = $wb.Worksheets.Add($wsTicketProd)
$detailProd.Name = "detail_prod"

# Header row
$detailProd.Range("B2").Value = "Name"
$detailProd.Range("C2").Value = "Data Type"
$detailProd.Range("D2").Value = "Primary Identifier"
$detailProd.Range("E2").Value = "Description"
$detailProd.Range("B2:D2").Font.Bold = $true
$detailProd.Range("B2:D2").Borders.LineStyle = 1
$detailProd.Range("B2:D2").Borders.Weight = 2
$detailProd.Range("E2").Font.Bold = $true
$detailProd.Range("E2").Borders.LineStyle = 1
$detailProd.Range("E2").Borders.Weight = 2
$detailProd.Range("E2").VerticalAlignment = -4108

# Row 3: dprod_id / Serial / TRUE / Llave primaria
$detailProd.Range("B3").Value = "dprod_id"
$detailProd.Range("C3").Value = "Serial"
$detailProd.Range("D3").Value = "TRUE"
$detailProd.Range("E3").Value = "Llave primaria"
$detailProd.Range("B3:D3").Borders.LineStyle = 1
$detailProd.Range("B3:D3").Borders.Weight = 2
$detailProd.Range("E3").Borders.LineStyle = 1
$detailProd.Range("E3").Borders.Weight = 2
$detailProd.Range("E3").VerticalAlignment = -4108

# Row 4: dtk_id / Integer / FALSE / id del detalle del ticket
$detailProd.Range("B4").Value = "dtk_id"
$detailProd.Range("C4").Value = "Integer"
$detailProd.Range("D4").Value = "FALSE"
$detailProd.Range("E4").Value = "id del detalle del ticket"
$detailProd.Range("B4:E4").Borders.LineStyle = 1
$detailProd.Range("B4:E4").Borders.Weight = 2
$detailProd.Range("B4:E4").VerticalAlignment = -4108

# Row 5: tpd_id / Integer / FALSE / id del ticket de produccion
$detailProd.Range("B5").Value = "tpd_id"
$detailProd.Range("C5").Value = "Integer"
$detailProd.Range("D5").Value = "FALSE"
$detailProd.Range("E5").Value = "id del ticket de produccion"
$detailProd.Range("B5:E5").Borders.LineStyle = 1
$detailProd.Range("B5:E5").Borders.Weight = 2

# Row 6: leftover formatted empty cell in column E
$detailProd.Range("E6").Font.Underline = 1

$detailProd.Columns.Item(4).ColumnWidth = 16.85546875
$detailProd.Columns.Item(5).ColumnWidth = 24.7109375

Write-Host "detail_prod created"

# --- Set active sheet / tab selection ---
$detailProd.Select()
$excel.ActiveWindow.ScrollWorkbookTabs(0)

Write-Host "Sheets now:"
foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
